$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.15635359287262
$ws.Range("B1").Value = 2.39827561378479
$ws.Range("D1").Value = 2.384132862091064
$ws.Range("E1").Value = 1.225741147994995
